$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('AA15').Value = 1
$ws.Range('AB15').Value = -1
$ws.Range('AC15').Value = 0.825
$ws.Range('B15').Value = 6992554
$ws.Range('F15').Value = 'Sukhothai FC'
$ws.Range('G15').Value = 'Trat FC'
$ws.Range('H15').Value = 0
$ws.Range('J15').Value = 'D'
$ws.Range('K15').Value = 1.8
$ws.Range('L15').Value = 3.6
$ws.Range('M15').Value = 4.333
$ws.Range('N15').Value = 1.833
$ws.Range('O15').Value = 3.75
$ws.Range('P15').Value = 4
$ws.Range('Q15').Value = -0.5
$ws.Range('R15').Value = 1.8
$ws.Range('S15').Value = 2
$ws.Range('T15').Value = 2.75
$ws.Range('W15').Value = -1
$ws.Range('X15').Value = 2.75
$ws.Range('Z15').Value = -1
$ws.Range('AA16').Value = -1
$ws.Range('AB16').Value = 0
$ws.Range('AC16').Value = 0
$ws.Range('B16').Value = 6992550
$ws.Range('F16').Value = 'Buriram United'
$ws.Range('G16').Value = 'Lamphun Warrior FC'
$ws.Range('H16').Value = 3
$ws.Range('J16').Value = 'H'
$ws.Range('K16').Value = 1.166
$ws.Range('L16').Value = 8
$ws.Range('M16').Value = 12
$ws.Range('N16').Value = 1.25
$ws.Range('O16').Value = 6.5
$ws.Range('P16').Value = 8.5
$ws.Range('Q16').Value = -1.75
$ws.Range('R16').Value = 1.875
$ws.Range('S16').Value = 1.925
$ws.Range('T16').Value = 3
$ws.Range('W16').Value = 0.25
$ws.Range('X16').Value = -1
$ws.Range('Z16').Value = 0.875
$ws.Range('AA117').Value = -1
$ws.Range('AB117').Value = -1
$ws.Range('AC117').Value = 0.95
$ws.Range('B117').Value = 7329293
$ws.Range('F117').Value = 'Chonburi'
$ws.Range('G117').Value = 'Bangkok United'
$ws.Range('H117').Value = 0
$ws.Range('I117').Value = 0
$ws.Range('K117').Value = 3.6
$ws.Range('L117').Value = 3.5
$ws.Range('M117').Value = 1.85
$ws.Range('N117').Value = 4.5
$ws.Range('O117').Value = 4
$ws.Range('P117').Value = 1.615
$ws.Range('Q117').Value = 0.75
$ws.Range('R117').Value = 1.975
$ws.Range('S117').Value = 1.825
$ws.Range('U117').Value = 1.85
$ws.Range('V117').Value = 1.95
$ws.Range('X117').Value = 3
$ws.Range('Z117').Value = 0.9750000000000001
$ws.Range('AA118').Value = 0.95
$ws.Range('AB118').Value = 0.825
$ws.Range('AC118').Value = -1
$ws.Range('B118').Value = 7485127
$ws.Range('F118').Value = 'BG Pathum United'
$ws.Range('G118').Value = 'Chiangrai Utd'
$ws.Range('H118').Value = 2
$ws.Range('I118').Value = 2
$ws.Range('K118').Value = 1.5
$ws.Range('L118').Value = 4
$ws.Range('M118').Value = 5.75
$ws.Range('N118').Value = 1.363
$ws.Range('O118').Value = 4.5
$ws.Range('P118').Value = 6.5
$ws.Range('Q118').Value = -1.25
$ws.Range('R118').Value = 1.85
$ws.Range('S118').Value = 1.95
$ws.Range('U118').Value = 1.825
$ws.Range('V118').Value = 1.975
$ws.Range('X118').Value = 3.5
$ws.Range('Z118').Value = -1
$ws.Range('B183').Value = 6992337
$ws.Range('E183').Value = 45389.3125
$ws.Range('F183').Value = 'Lamphun Warrior FC'
$ws.Range('G183').Value = 'Police Tero FC'
$ws.Range('K183').Value = 1.444
$ws.Range('L183').Value = 4.2
$ws.Range('M183').Value = 5.5
$ws.Range('N183').Value = 1.5
$ws.Range('O183').Value = 4.2
$ws.Range('P183').Value = 5.25
$ws.Range('Q183').Value = -1.25
$ws.Range('R183').Value = 1.975
$ws.Range('S183').Value = 1.825
$ws.Range('T183').Value = 3.25
$ws.Range('U183').Value = 1.975
$ws.Range('V183').Value = 1.825
$ws.Range('B184').Value = 6992704
$ws.Range('E184').Value = 45389.33333333334
$ws.Range('F184').Value = 'Bangkok United'
$ws.Range('G184').Value = 'Port FC'
$ws.Range('K184').Value = 1.85
$ws.Range('L184').Value = 3.6
$ws.Range('M184').Value = 3.4
$ws.Range('N184').Value = 1.727
$ws.Range('O184').Value = 3.75
$ws.Range('P184').Value = 3.8
$ws.Range('Q184').Value = -0.75
$ws.Range('R184').Value = 1.975
$ws.Range('S184').Value = 1.825
$ws.Range('T184').Value = 2.75
$ws.Range('U184').Value = 1.825
$ws.Range('V184').Value = 1.975
$ws.Range('B185').Value = 6992702
$ws.Range('E185').Value = 45389.35416666666
$ws.Range('F185').Value = 'Chonburi'
$ws.Range('G185').Value = 'Ratchaburi FC'
$ws.Range('K185').Value = 2
$ws.Range('L185').Value = 3.3
$ws.Range('M185').Value = 3.25
$ws.Range('N185').Value = 2.1
$ws.Range('O185').Value = 3.25
$ws.Range('P185').Value = 3
$ws.Range('Q185').Value = -0.25
$ws.Range('R185').Value = 1.9
$ws.Range('S185').Value = 1.9
$ws.Range('T185').Value = 2.75
$ws.Range('U185').Value = 1.9
$ws.Range('V185').Value = 1.9
$ws.Range('B186').Value = 6992708
$ws.Range('E186').Value = 45389.375
$ws.Range('F186').Value = 'Sukhothai FC'
$ws.Range('G186').Value = 'Nakhon Pathom FC'
$ws.Range('K186').Value = 1.85
$ws.Range('L186').Value = 3.4
$ws.Range('M186').Value = 3.6
$ws.Range('N186').Value = 1.85
$ws.Range('O186').Value = 3.4
$ws.Range('P186').Value = 3.6
$ws.Range('Q186').Value = -0.5
$ws.Range('R186').Value = 1.9
$ws.Range('S186').Value = 1.9
$ws.Range('U186').Value = 1.95
$ws.Range('V186').Value = 1.85
